$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# [FEATURE] Migracion UFT - ANSES, Bloq y Desbloq
# ANS01, ANS02
# BYD21 al BYD26

# Row 57 - ANS01
$ws.Range("A57").Value = "F00020"
$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").HorizontalAlignment = -4152
$ws.Range("C57").Value = "020"

# Row 58 - ANS02
$ws.Range("A58").Value = "F00847"
$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").HorizontalAlignment = -4152
$ws.Range("C58").Value = "074"

# Row 59
$ws.Range("A59").Value = "F03808"
$ws.Range("C59").Value = 322

$ws.Range("C60").Select()
